$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.770.92"
$ws.Range("E2").Value = "  +0.41%  "
$ws.Range("D3").Value = "2.101.89"
$ws.Range("E3").Value = "  +0.21%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "228.69"
$ws.Range("E5").Value = "  -0.18%  "
$ws.Range("D6").Value = "0.617"
$ws.Range("E6").Value = "  +0.18%  "
$ws.Range("D7").Value = "62.48"
$ws.Range("E7").Value = "  +2.02%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").Value = "0.391"
$ws.Range("E9").Value = "  +2.26%  "
$ws.Range("D10").Value = "0.0841"
$ws.Range("E10").Value = "  -0.19%  "
$ws.Range("D11").Value = "0.103"
$ws.Range("E11").Value = "  -1.00%  "
$ws.Range("D12").Value = "15.73"
$ws.Range("E12").Value = "  +6.38%  "
$ws.Range("D13").Value = "2.413.59"
$ws.Range("E13").Value = "  +0.20%  "
$ws.Range("E14").Value = "  -1.69%  "
$ws.Range("D15").Value = "0.811"
$ws.Range("E15").Value = "  +3.76%  "
$ws.Range("E16").Value = "  +0.63%  "
$ws.Range("D17").Value = "2.103.86"
$ws.Range("E17").Value = "  -3.33%  "
$ws.Range("D18").Value = "38.809.13"
$ws.Range("E18").Value = "  +0.78%  "
$ws.Range("D19").Value = "71.95"
$ws.Range("E19").Value = "  +1.62%  "
$ws.Range("D20").Value = "6.10"
$ws.Range("E20").Value = "  +1.34%  "
$ws.Range("E21").Value = "  +0.64%  "
$ws.Range("D22").Value = "228.78"
$ws.Range("E22").Value = "  +0.87%  "
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("E24").Value = "  -3.77%  "
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("D26").Value = "171.84"
$ws.Range("D27").Value = "9.58"
$ws.Range("E27").Value = "  +1.73%  "
$ws.Range("E28").Value = "  +5.87%  "
$ws.Range("E29").Value = "  +4.75%  "
$ws.Range("E30").Value = "  +1.52%  "
$ws.Range("E31").Value = "  +8.46%  "
$ws.Range("E32").Value = "  +0.66%  "
$ws.Range("E33").Value = "  +1.80%  "
$ws.Range("E34").Value = "  -0.27%  "
$ws.Range("E35").Value = "  +7.05%  "
$ws.Range("D36").Value = "0.0618"
$ws.Range("E36").Value = "  +2.01%  "
$ws.Range("E37").Value = "  +0.99%  "
$ws.Range("D38").Value = "3.59"
$ws.Range("E38").Value = "  +1.19%  "
$ws.Range("D39").Value = "0.999"
$ws.Range("E39").Value = "  -0.05%  "
$ws.Range("D40").Value = "18.01"
$ws.Range("E40").Value = "  -2.64%  "
$ws.Range("D41").Value = "102.97"
$ws.Range("E41").Value = "  +3.05%  "
$ws.Range("D42").Value = "0.0229"
$ws.Range("E42").Value = "  +3.55%  "
$ws.Range("D43").Value = "1.533.06"
$ws.Range("E43").Value = "  -0.77%  "
$ws.Range("E44").Value = "  +4.12%  "
$ws.Range("E45").Value = "  +4.22%  "
$ws.Range("E46").Value = "  -1.24%  "
$ws.Range("D47").Value = "0.0910"
$ws.Range("E47").Value = "  -0.37%  "
$ws.Range("E48").Value = "  -2.20%  "
$ws.Range("E49").Value = "  +1.22%  "
$ws.Range("E50").Value = "  -0.18%  "
$ws.Range("D51").Value = "2.299.98"
$ws.Range("E51").Value = "  +0.16%  "
